$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (s="1") from E1 to F1, then set the header text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:40:48.031806",
    "2021-10-05 13:40:48.031818",
    "2021-10-05 13:40:48.031822",
    "2021-10-05 13:40:48.031825",
    "2021-10-05 13:40:48.031828",
    "2021-10-05 13:40:48.031832",
    "2021-10-05 13:40:48.031835",
    "2021-10-05 13:40:48.031838",
    "2021-10-05 13:40:48.031841",
    "2021-10-05 13:40:48.031844",
    "2021-10-05 13:40:48.031847",
    "2021-10-05 13:40:48.031850",
    "2021-10-05 13:40:48.031853"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
